$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '36.555.40'; 'E2' = '  -0.43%  '
    'D3' = '1.967.48'; 'E3' = '  +0.15%  '
    'E4' = '  +0.01%  '
    'D5' = '243.72'; 'E5' = '  -0.34%  '
    'D6' = '0.619'; 'E6' = '  +0.63%  '
    'D7' = '59.94'; 'E7' = '  +2.47%  '
    'E8' = '  -0.02%  '
    'D9' = '0.378'; 'E9' = '  +1.97%  '
    'E10' = '  -3.01%  '
    'E11' = '  +0.63%  '
    'D12' = '14.22'; 'E12' = '  +3.75%  '
    'E13' = '  +2.37%  '
    'D14' = '2.257.12'; 'E14' = '  +0.23%  '
    'D15' = '21.61'; 'E15' = '  -2.84%  '
    'E16' = '  +0.37%  '
    'D17' = '1.967.58'; 'E17' = '  +0.06%  '
    'D18' = '36.519.59'; 'E18' = '  -0.50%  '
    'E19' = '  +0.06%  '
    'D20' = '0.0₃0853'; 'E20' = '  -1.10%  '
    'D21' = '229.36'; 'E21' = '  +0.49%  '
    'E22' = '  -0.44%  '
    'D23' = '0.999'; 'E23' = '  -0.11%  '
    'E24' = '  +1.35%  '
    'E25' = '  +1.03%  '
    'E26' = '  +5.31%  '
    'D27' = '9.11'; 'E27' = '  -2.27%  '
    'D28' = '162.27'; 'E28' = '  +0.83%  '
    'D29' = '19.38'; 'E29' = '  +0.02%  '
    'E30' = '  +19.55%  '
    'E31' = '  +2.12%  '
    'D32' = '4.79'; 'E32' = '  +2.49%  '
    'E33' = '  -0.76%  '
    'E34' = '  +6.79%  '
    'E35' = '  +3.78%  '
    'E36' = '  -0.08%  '
    'D37' = '3.34'; 'E37' = '  -2.93%  '
    'E38' = '  +0.33%  '
    'D39' = '5.39'; 'E39' = '  -14.14%  '
    'D40' = '0.0965'; 'E40' = '  -3.61%  '
    'E41' = '  +0.49%  '
    'E42' = '  +0.17%  '
    'E44' = '  -1.20%  '
    'D45' = '1.365.37'; 'E45' = '  +1.40%  '
    'D46' = '88.97'; 'E46' = '  +1.65%  '
    'E47' = '  -1.61%  '
    'E48' = '  +0.86%  '
    'D49' = '2.82'; 'E49' = '  -0.86%  '
    'D50' = '45.87'; 'E50' = '  +5.30%  '
    'D51' = '2.151.07'; 'E51' = '  +0.37%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
